# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record is inserted at row 41 (pushing the existing
# rows 41-164 down to 42-165), extending the used range to A1:R165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 41; this shifts rows 41..164
# down to 42..165 (values, styles and the sheet dimension all move with
# it), matching the diff exactly.
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new weekly record.
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(41, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(41, 3).Value = "Los Lagos"
$ws.Cells.Item(41, 4).Value = 44525
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(41, 6).Value = 100112017
$ws.Cells.Item(41, 7).Value = "Apio"
$ws.Cells.Item(41, 8).Value = "Americana (o)"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 35
$ws.Cells.Item(41, 11).Value = 12000
$ws.Cells.Item(41, 12).Value = 12000
$ws.Cells.Item(41, 13).Value = 12000
$ws.Cells.Item(41, 14).Value = "$/docena de matas"
$ws.Cells.Item(41, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(41, 16).Value = 2000
$ws.Cells.Item(41, 17).Value = 6
$ws.Cells.Item(41, 18).Value = "Hortaliza"
